$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5,6,7 hold Malik Monk / Keegan Murray / Kevin Huerter.
# The roster is being re-sorted so Keegan Murray moves to row 5,
# Kevin Huerter moves to row 6, and Malik Monk moves to row 7
# (column A, the roster index, stays put).

$ws.Range("B5").Value = 13
$ws.Range("C5").Value = "Keegan Murray"
$ws.Range("D5").Value = "SF"
$ws.Range("E5").Value = "6-8"
$ws.Range("F5").Value = 215
$ws.Range("G5").Value = "August 19, 2000"
$ws.Range("I5").Value = "R"
$ws.Range("J5").Value = "Iowa"
$ws.Range("K5").Value = "https://www.basketball-reference.com/players/m/murrake02.html"

$ws.Range("B6").Value = 9
$ws.Range("C6").Value = "Kevin Huerter"
$ws.Range("D6").Value = "SG"
$ws.Range("E6").Value = "6-7"
$ws.Range("F6").Value = 190
$ws.Range("G6").Value = "August 27, 1998"
$ws.Range("I6").Value = "4"
$ws.Range("J6").Value = "Maryland"
$ws.Range("K6").Value = "https://www.basketball-reference.com/players/h/huertke01.html"

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = "Malik Monk"
$ws.Range("D7").Value = "SG"
$ws.Range("E7").Value = "6-3"
$ws.Range("F7").Value = 200
$ws.Range("G7").Value = "February 4, 1998"
$ws.Range("I7").Value = "5"
$ws.Range("J7").Value = "Kentucky"
$ws.Range("K7").Value = "https://www.basketball-reference.com/players/m/monkma01.html"

# Kessler Edwards (row 18) gets a jersey number assigned: No. = 17
$ws.Range("B18").Value = 17
